$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Name"), shifting Name/Phone number/Party
# size one column to the right, to make room for the new "Deleted" column.
$ws.Columns("B").Insert()

# Add the new "Deleted_ý" header and rename header "Id" -> "Id_ý" (in this
# order, so the shared-strings table gets the two new strings in the same
# sequence Excel produced them).
$ws.Range("B1").Value = "Deleted_ý"
$ws.Range("A1").Value = "Id_ý"

# Give the new column the same width as column A (best effort - keeps the
# original custom width flag on the inserted column).
$ws.Columns("B").ColumnWidth = 13.5

# Populate the new "Deleted" column with boolean FALSE for every data row.
$ws.Range("B2").Value = $false
$ws.Range("B2").Copy()
$ws.Range("B3:B4").PasteSpecial(-4122)
$ws.Range("B3").Value = $false
$ws.Range("B4").Value = $false

# Rename the built-in "Normal" cell style to the Dutch localisation "Standaard".
$wb.Styles.Item("Normal").Name = "Standaard"

# Match the final selection left behind in the saved file.
$ws.Range("B3").Select()
